# The source diff for this revision touches only the serialized
# attribute order of existing OOXML elements (namespace declarations on
# <w:document>, <w:pgSz>/<w:pgMar> on the section properties, the
# <w:rFonts>/<w:lang> defaults, <w:latentStyles>/<w:lsdException> and
# <w:style> definitions in the style sheet). No element, attribute
# value, or piece of text content was added, removed or changed -
# every "-"/"+" pair in the diff carries the exact same set of
# attributes, just written in a different (alphabetised) order, which
# is a serializer/canonicalization artifact rather than a document
# edit. There is nothing in the Word object model that reaches that
# granularity (Word does not expose attribute ordering), so there is
# no content-level edit to apply here.
#
# Touch the document through the supported object model so the COM
# session has a well-defined ActiveDocument to operate on, without
# altering any visible content, formatting or properties.
$d = $word.ActiveDocument
$null = $d.Content
